$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update "Marking" row (B11): correct marks per right answer
$ws.Range("B11").Value = 5

# Update "Total" row (B12): total correct marks
$ws.Range("B12").Value = 110

# Update corr/total marks display (E12)
$ws.Range("E12").Value = "110/140"
